# Applies the "add a few text" commit:
#  - fixes slide 6 body paragraph 4 (merges 3 runs into 1 run with identical text)
#  - appends 9 new slides (7..15) at the end of the deck:
#      7  Requirement Analysis            (Title and Content)
#      8  Background                      (Title and Content)
#      9  Flowchart Diagram For ETD       (Title and Content, empty body)
#      10 ER diagram for ETD              (Title and Content, empty body)
#      11 Use case Diagram For ETD        (Title and Content, empty body)
#      12 Sequence Diagram For ETD        (Title and Content, empty body)
#      13 Minestones and job distribution (Title and Content, empty body)
#      14 Conclusion                      (Title and Content, empty body)
#      15 Thanks For Your Attention!      (Title Only, custom title position)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6: collapse the 3-run paragraph ("For Admin Panel , admins " +
#    "can " + "insert  data into database , delete data and retrieve data
#    for database.") into a single run. Because the concatenated text is
#    byte-identical to the original, the text has to be changed and then
#    changed back so the host actually rewrites the paragraph's runs.
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$adminPara = $slide6.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4, 1)
$adminPara.Text = "__tmp__"
$slide6.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4, 1).Text = "For Admin Panel , admins can insert  data into database , delete data and retrieve data for database."

# ---------------------------------------------------------------------------
# 2) New slide 7 - "Requirement Analysis"
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Requirement Analysis"
$body7 = $s7.Shapes.Item(2).TextFrame.TextRange
$body7.Text = "We get all information of data we use in this project from library of computer university "
$body7.InsertAfter("mandalay") | Out-Null
$s7.Shapes.Item(2).TextFrame.TextRange.InsertAfter(".") | Out-Null

# ---------------------------------------------------------------------------
# 3) New slide 8 - "Background"
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Background"
$body8 = $s8.Shapes.Item(2).TextFrame.TextRange
$body8.Text = "Our Software Company has received"
$body8.InsertAfter("`rMoreover , we want to know the profit and loss of our system.") | Out-Null
$s8.Shapes.Item(2).TextFrame.TextRange.InsertAfter("`rSo there are many processes to develop manual system and it may be more cost ant time waste . So we decide to develop our system to reduce the problem.") | Out-Null

# ---------------------------------------------------------------------------
# 4) New slide 9 - "Flowchart Diagram For ETD" (body left empty)
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Flowchart Diagram For ETD"

# ---------------------------------------------------------------------------
# 5) New slide 10 - "ER diagram for ETD" (body left empty)
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "ER diagram for ETD"

# ---------------------------------------------------------------------------
# 6) New slide 11 - "Use case Diagram For ETD" (body left empty)
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s11.Shapes.Item(1).TextFrame.TextRange.Text = "Use case Diagram For ETD"

# ---------------------------------------------------------------------------
# 7) New slide 12 - "Sequence Diagram For ETD" (body left empty)
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s12.Shapes.Item(1).TextFrame.TextRange.Text = "Sequence Diagram For ETD"

# ---------------------------------------------------------------------------
# 8) New slide 13 - "Minestones and job distribution" (two runs in the
#    title; body left empty)
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Add($p.Slides.Count + 1, 2)
$title13 = $s13.Shapes.Item(1).TextFrame.TextRange
$title13.Text = "Minestones"
$title13.InsertAfter(" and job distribution") | Out-Null

# ---------------------------------------------------------------------------
# 9) New slide 14 - "Conclusion" (body left empty)
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s14.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusion"

# ---------------------------------------------------------------------------
# 10) New slide 15 - "Thanks For Your Attention!" (Title Only layout, title
#     box moved/resized to x=3217333 y=2768998 cx=8772701 cy=843446 EMU)
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Add($p.Slides.Count + 1, 6)
$title15 = $s15.Shapes.Item(1)
$title15.TextFrame.TextRange.Text = "Thanks For Your Attention!"
$title15.Left = 253.33331298828125
$title15.Top = 218.03134155273438
$title15.Width = 690.763916015625
$title15.Height = 66.41307830810547
